$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...(2) a projected" -> split "projected" into its own run
#    ("DITCO) with: ... (2) a " | "projected")
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("nature; (2) a projected", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Start = $rng1.End - 9   # length of "projected"
    # Forcing an explicit (but value-preserving) font write makes Word break
    # the run in two without altering the visible formatting: flip the color
    # away and back so the engine mints a fresh run/rPr for this sub-range.
    $rng1.Font.Color = 255
    $rng1.Font.Color = 0
}

# ---------------------------------------------------------------------------
# 2) "nse Federal Acquisition Regulation Supplement " -> split into
#    "nse Federal " | "Acquisition Regulation Supplement "
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Acquisition Regulation Supplement ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Bold = 1
    $rng2.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) "PART 1 - CERTIFICATION OF SEVERABILITY  " -> split into
#    "PART 1 - CERTIFICA" | "TION OF SEVERABILITY  "
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("TION OF SEVERABILITY  ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Font.Color = 255
    $rng3.Font.Color = 0
}

# ---------------------------------------------------------------------------
# 4) second "{missionOwner}" (the one under the PM/PM/COR signature block,
#    just before the Financial POC block) -> "{primaryContact}"
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "{missionOwner}") {
        $pr = $p.Range
        $pr.Find.Execute("missionOwner", $false, $false, $false, $false, $false, $true, 1, $false, "primaryContact", 2)
    }
}
